# Updated cryptos list on Wed Aug 28 13:57:35 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "59.520.54"
$ws.Range("E2").Value = "  -4.47%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.480.63"
$ws.Range("E3").Value = "  -5.06%  "

# Row 4
$ws.Range("E4").Value = "  +0.13%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "539.87"
$ws.Range("E5").Value = "  -2.87%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.39"
$ws.Range("E6").Value = "  -5.49%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.997"
$ws.Range("E7").Value = "  -0.24%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.574"
$ws.Range("E8").Value = "  -3.22%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.507.75"
$ws.Range("E9").Value = "  -4.10%  "

# Row 10
$ws.Range("E10").Value = "  -3.64%  "

# Row 11
$ws.Range("E11").Value = "  -1.37%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.33"
$ws.Range("E12").Value = "  -2.58%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.357"
$ws.Range("E13").Value = "  -2.14%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.924.58"
$ws.Range("E14").Value = "  -4.78%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "24.25"
$ws.Range("E15").Value = "  -5.70%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "59.378.54"
$ws.Range("E16").Value = "  -4.49%  "

# Row 17
$ws.Range("E17").Value = "  -3.34%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.503.84"
$ws.Range("E18").Value = "  -4.11%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.45"
$ws.Range("E19").Value = "  -1.43%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.36"
$ws.Range("E20").Value = "  -3.86%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "326.13"
$ws.Range("E21").Value = "  -4.42%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.996"
$ws.Range("E22").Value = "  -0.37%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.77"
$ws.Range("E23").Value = "  -5.17%  "

# Row 24
$ws.Range("E24").Value = "  -2.64%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.445"
$ws.Range("E25").Value = "  -10.77%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.162"
$ws.Range("E26").Value = "  -3.36%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.994"
$ws.Range("E27").Value = "  -0.65%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.607.33"
$ws.Range("E28").Value = "  -4.55%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.83"
$ws.Range("E29").Value = "  -2.50%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.10"
$ws.Range("E30").Value = "  -0.48%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.29"
$ws.Range("E31").Value = "  -3.49%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0₃0783"
$ws.Range("E32").Value = "  -5.70%  "

# Row 33
$ws.Range("E33").Value = "  -4.45%  "

# Row 34
$ws.Range("B34").Value = "USDe"
$ws.Range("C34").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.997"
$ws.Range("E34").Value = "  -0.22%  "

# Row 35
$ws.Range("B35").Value = "Monero"
$ws.Range("C35").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "159.46"
$ws.Range("E35").Value = "  -0.65%  "

# Row 36
$ws.Range("E36").Value = "  +0.62%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "18.67"
$ws.Range("E37").Value = "  -3.11%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.47"
$ws.Range("E38").Value = "  -4.89%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.66"
$ws.Range("E39").Value = "  -4.42%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.96"
$ws.Range("E40").Value = "  -2.28%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "314.17"
$ws.Range("E41").Value = "  -7.24%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "36.68"
$ws.Range("E42").Value = "  -2.65%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.75"
$ws.Range("E43").Value = "  -3.64%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.832"
$ws.Range("E44").Value = "  -6.86%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.995"
$ws.Range("E45").Value = "  -0.33%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.599"
$ws.Range("E46").Value = "  -1.76%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.78"
$ws.Range("E47").Value = "  -1.72%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "125.19"
$ws.Range("E48").Value = "  -0.89%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0937"
$ws.Range("E49").Value = "  -3.01%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0526"
$ws.Range("E50").Value = "  -4.01%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0230"
$ws.Range("E51").Value = "  -3.97%  "

